$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.267.12'
$ws.Range("E2").Value = '  +1.38%  '

$ws.Range("D3").Value = '1.807.65'
$ws.Range("E3").Value = '  +3.15%  '

$ws.Range("E4").Value = '  -0.30%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '336.60'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.30%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9985'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.26%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4606'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +20.66%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3716'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +9.64%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.25'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.46%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.156'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.07%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07657'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.25%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.47'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.39%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.0000'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.31%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.356'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.41%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.478'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.91%  '

$ws.Range("D16").Value = '1.804.72'
$ws.Range("E16").Value = '  +2.66%  '

$ws.Range("E17").Value = '  +4.14%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06714'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.54%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '81.98'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.79%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9988'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.17%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.52'
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.444'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.75%  '

$ws.Range("D23").Value = '28.270.40'
$ws.Range("E23").Value = '  +1.33%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.92'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.88%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.417'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.36%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.86'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.49%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '153.60'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.00%  '

$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.391'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.55%  '

$ws.Range("D29").Value = '2.010.22'
$ws.Range("E29").Value = '  +2.65%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.78'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.49%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.261'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.80%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.032'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.24%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09570'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +8.81%  '

$ws.Range("E34").Value = '  +1.30%  '

$ws.Range("E35").Value = '  +6.03%  '

$ws.Range("B36").Value = 'Aptos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.17'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.05%  '

$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06372'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.40%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02362'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.55%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.285'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.36%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6664'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.43%  '

$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.241'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.03%  '

$ws.Range("B42").Value = 'WEMIXTOKEN'
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.500'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.72%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.163'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.49%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.30'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.03%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9980'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.26%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6139'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.58%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.829'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.17%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '130.03'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.17%  '

$ws.Range("E49").Value = '  +3.04%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07164'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.61%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.177'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.97%  '
